$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 94: correct the intraday timestamp in column A (date) ---
# Old value 45477.6495949074 (2024-07-04 15:35:25) -> new value
# 45477.2916666667 (2024-07-04 07:00:00), matching the rest of the
# column's "07:00:00" time-of-day convention.
$ws.Range("A94").Value = 45477.2916666667

# --- New row 95: append the next day's OHLCV record ---
# Copy row 94's date cell (A94) into A95 first so the new date cell
# inherits the same date/time number format (style index 1), then
# overwrite its value with the new timestamp.
$ws.Range("A94").Copy()
$ws.Range("A95").PasteSpecial(-4122)
$ws.Range("A95").Value = 45478.4268518519

$ws.Range("B95").Value = 37500
$ws.Range("C95").Value = 3.46000003814697
$ws.Range("D95").Value = 3.29999995231628
$ws.Range("E95").Value = 3.4300000667572
$ws.Range("F95").Value = 3.44000005722046

# Column G ("adj_close") stores the same number as column F, but as
# text (shared string) rather than a numeric cell, matching every
# other row in the sheet. A leading apostrophe forces text storage of
# a numeric-looking literal; resetting the style back to "Normal"
# avoids leaving a stray text-number-format style behind.
$ws.Range("G95").Value = "'3.44000005722046"
$ws.Range("G95").Style = "Normal"

# Column H ("ticker") is always "ESPE.MI".
$ws.Range("H95").Value = "ESPE.MI"
